$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row being appended to the data table (CompaNanny Statenkwartier BSO, 2024-09-17)
$row = 83

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Statenkwartier BSO"
$ws.Cells.Item($row, 3).Value = "VGO"

# Leading apostrophe forces the date-like string to be stored as text
# (matching the other rows, which hold plain "yyyy-mm-dd" text, not a date).
$ws.Cells.Item($row, 4).Value = "'2024-09-17"
# Re-apply the plain default style used by the rest of the data rows, so the
# quote-prefix marker doesn't leave a stray format on the cell.
$ws.Cells.Item($row, 4).Style = $ws.Cells.Item($row - 1, 4).Style

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
